$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 'flower/flower078.png'
$ws.Range("D2").Value = 'kehren'
$ws.Range("E2").Value = 'flower'

$ws.Range("B3").Value = 119
$ws.Range("C3").Value = 'dog/dog119.png'
$ws.Range("D3").Value = 'fesseln'
$ws.Range("E3").Value = 'dog'

$ws.Range("B4").Value = 63
$ws.Range("C4").Value = 'dog/dog099.png'
$ws.Range("D4").Value = 'füttern'
$ws.Range("E4").Value = 'dog'

$ws.Range("B5").Value = 76
$ws.Range("C5").Value = 'dog/dog097.png'
$ws.Range("D5").Value = 'spielen'
$ws.Range("E5").Value = 'dog'

$ws.Range("B6").Value = 114
$ws.Range("C6").Value = 'dog/dog086.png'
$ws.Range("D6").Value = 'sondern'
$ws.Range("E6").Value = 'dog'

$ws.Range("B7").Value = 74
$ws.Range("C7").Value = 'flower/flower077.png'
$ws.Range("D7").Value = 'formen'
$ws.Range("E7").Value = 'flower'

$ws.Range("B8").Value = 122
$ws.Range("C8").Value = 'dog/dog121.png'
$ws.Range("D8").Value = 'kaufen'
$ws.Range("E8").Value = 'dog'

$ws.Range("B9").Value = 31
$ws.Range("C9").Value = 'dog/dog073.png'
$ws.Range("D9").Value = 'fliehen'
$ws.Range("E9").Value = 'dog'

$ws.Range("B10").Value = 62
$ws.Range("C10").Value = 'flower/flower088.png'
$ws.Range("D10").Value = 'schicken'
$ws.Range("E10").Value = 'flower'

$ws.Range("B11").Value = 85
$ws.Range("C11").Value = 'flower/flower074.png'
$ws.Range("D11").Value = 'pflegen'
$ws.Range("E11").Value = 'flower'

$ws.Range("B12").Value = 90
$ws.Range("C12").Value = 'dog/dog110.png'
$ws.Range("D12").Value = 'rasen'
$ws.Range("E12").Value = 'dog'

$ws.Range("B13").Value = 61
$ws.Range("C13").Value = 'dog/dog106.png'
$ws.Range("D13").Value = 'fühlen'
$ws.Range("E13").Value = 'dog'

$ws.Range("B14").Value = 42
$ws.Range("C14").Value = 'flower/flower098.png'
$ws.Range("D14").Value = 'wiegen'
$ws.Range("E14").Value = 'flower'

$ws.Range("B15").Value = 127
$ws.Range("C15").Value = 'dog/dog122.png'
$ws.Range("D15").Value = 'haken'
$ws.Range("E15").Value = 'dog'

$ws.Range("B16").Value = 39
$ws.Range("C16").Value = 'flower/flower097.png'
$ws.Range("D16").Value = 'gelten'
$ws.Range("E16").Value = 'flower'

$ws.Range("B17").Value = 11
$ws.Range("C17").Value = 'flower/flower089.png'
$ws.Range("D17").Value = 'biegen'
$ws.Range("E17").Value = 'flower'

$ws.Range("B18").Value = 53
$ws.Range("C18").Value = 'dog/dog068.png'
$ws.Range("D18").Value = 'nehmen'
$ws.Range("E18").Value = 'dog'

$ws.Range("B19").Value = 106
$ws.Range("C19").Value = 'flower/flower124.png'
$ws.Range("D19").Value = 'tagen'
$ws.Range("E19").Value = 'flower'

$ws.Range("B20").Value = 69
$ws.Range("C20").Value = 'dog/dog066.png'
$ws.Range("D20").Value = 'drehen'
$ws.Range("E20").Value = 'dog'

$ws.Range("B21").Value = 9
$ws.Range("C21").Value = 'dog/dog065.png'
$ws.Range("D21").Value = 'runden'
$ws.Range("E21").Value = 'dog'

$ws.Range("B22").Value = 46
$ws.Range("C22").Value = 'flower/flower101.png'
$ws.Range("D22").Value = 'segeln'
$ws.Range("E22").Value = 'flower'

$ws.Range("B23").Value = 116
$ws.Range("C23").Value = 'flower/flower115.png'
$ws.Range("D23").Value = 'tauschen'
$ws.Range("E23").Value = 'flower'

$ws.Range("B24").Value = 13
$ws.Range("C24").Value = 'flower/flower094.png'
$ws.Range("D24").Value = 'fliegen'
$ws.Range("E24").Value = 'flower'

$ws.Range("B25").Value = 40
$ws.Range("C25").Value = 'flower/flower103.png'
$ws.Range("D25").Value = 'posten'
$ws.Range("E25").Value = 'flower'

$ws.Range("B26").Value = 81
$ws.Range("C26").Value = 'dog/dog108.png'
$ws.Range("D26").Value = 'enden'
$ws.Range("E26").Value = 'dog'

$ws.Range("B27").Value = 72
$ws.Range("C27").Value = 'dog/dog083.png'
$ws.Range("D27").Value = 'scheitern'
$ws.Range("E27").Value = 'dog'

$ws.Range("B28").Value = 22
$ws.Range("C28").Value = 'flower/flower081.png'
$ws.Range("D28").Value = 'opfern'
$ws.Range("E28").Value = 'flower'

$ws.Range("B29").Value = 30
$ws.Range("C29").Value = 'dog/dog080.png'
$ws.Range("D29").Value = 'stechen'
$ws.Range("E29").Value = 'dog'

$ws.Range("B30").Value = 25
$ws.Range("C30").Value = 'flower/flower100.png'
$ws.Range("D30").Value = 'laufen'
$ws.Range("E30").Value = 'flower'

$ws.Range("B31").Value = 34
$ws.Range("C31").Value = 'flower/flower085.png'
$ws.Range("D31").Value = 'loben'
$ws.Range("E31").Value = 'flower'

$ws.Range("B32").Value = 23
$ws.Range("C32").Value = 'dog/dog095.png'
$ws.Range("D32").Value = 'währen'
$ws.Range("E32").Value = 'dog'

$ws.Range("B33").Value = 125
$ws.Range("C33").Value = 'flower/flower117.png'
$ws.Range("D33").Value = 'ehren'
$ws.Range("E33").Value = 'flower'
